# regen sval data to filter save games
# Updates columns B:G (TB, d2S, K, IP, Win, sum) for rows 2-15 on the active sheet
# with recomputed stat values (Win column F is unchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1.455362044514542, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 0, 3.754798637575387),
    @(0.01293466051926884, 0.306821227259698, 0.1494219747398047, 0.4942365360607697, 0, 0.9634143985795411),
    @(0.6606524410359556, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 1, 2.960089034096801),
    @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 1, 6.189590430959694),
    @(0.1190320826869504, 0.04071648406533734, 0.7527432677738641, 0.4942365360607697, 0, 1.406728370586922),
    @(1.455362044514542, 0.306821227259698, 0.7527432677738641, 0.4942365360607697, 1, 3.009163075608874),
    @(3.286832544864788, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 1, 5.586269137925634),
    @(3.286832544864788, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 1, 8.974608811992548),
    @(0.2917716402565462, 0.306821227259698, 0.7527432677738641, 0.4942365360607697, 0, 1.845572671350878),
    @(3.286832544864788, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 1, 5.586269137925634),
    @(1.455362044514542, 0.306821227259698, 0.1494219747398047, 0.4942365360607697, 0, 2.405841782574814),
    @(1.455362044514542, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 1, 4.358119930609447),
    @(3.286832544864788, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 0, 5.586269137925634),
    @(3.286832544864788, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 1, 8.974608811992548)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
